# Replace the "2018 Campaign Dates that use Perseus: Oct. 30-Nov. 8 and
# Nov. 29-Dec. 8" paragraphs with the newly-translated Orion campaign
# dates. Each of the four occurrences in the document is collapsed down
# to a single, unformatted run containing the new sentence (any leading
# <w:br/> run some copies carry is removed along with the rest).

$d = $word.ActiveDocument
$newText = "Campaign Dates that use Orion: January 16-25, February 14-23, March 14-24"
$needle = "Campaign Dates that use Perseus"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*$needle*") {
        $r = $p.Range
        [void]$r.MoveEnd(1, -1)
        $r.Delete()
        $r.InsertAfter($newText)
    }
}
